# Applies the "languages support and all" change:
#  1. Insert a new centered paragraph for {{ personal_details.job_title }}
#     right after the name paragraph.
#  2. Append two new runs (DOB / Gender) to the contact-details paragraph,
#     right after the existing Address run.
#  3. Insert a new "Languages" section (if/heading/for-loop/endif) right
#     after the Skills section's closing {% endif %}.

$d = $word.ActiveDocument
$CR = [char]13

# ---------------------------------------------------------------------
# 1) New paragraph: {{ personal_details.job_title }}, centered, directly
#    under the {{ personal_details.name }} title paragraph.
# ---------------------------------------------------------------------
$nameRange = $d.Paragraphs.Item(1).Range
$nameRange.Collapse(0)
$nameRange.InsertParagraphAfter()

$jobTitlePara = $d.Paragraphs.Item(2)
$jobTitlePara.Style = $d.Styles.Item("Normal")
$jobTitlePara.Alignment = 1
$jobTitlePara.Range.InsertAfter("{{ personal_details.job_title }}")

# ---------------------------------------------------------------------
# 2) Two new runs appended after the Address run in the contact-details
#    paragraph (now paragraph 3, since we inserted a paragraph above).
# ---------------------------------------------------------------------
$contactPara = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($contactPara.Range.End - 1, $contactPara.Range.End - 1)
$insertPoint.InsertAfter("{% if personal_details.date_of_birth %} | DOB: {{ personal_details.date_of_birth }}{% endif %}")

$contactPara = $d.Paragraphs.Item(3)
$insertPoint = $d.Range($contactPara.Range.End - 1, $contactPara.Range.End - 1)
$insertPoint.InsertAfter("{% if personal_details.gender %} | Gender: {{ personal_details.gender }}{% endif %}")

# ---------------------------------------------------------------------
# 3) New "Languages" section, inserted right after the Skills section's
#    closing {% endif %} paragraph (i.e. right before {% if education %}).
# ---------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd($CR)
    if ($txt -eq "{% endif %}") {
        $nextP = $d.Paragraphs.Item($i + 1)
        $nextTxt = $nextP.Range.Text.TrimEnd($CR)
        if ($nextTxt -eq "{% if education %}") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the Skills section's closing {% endif %} paragraph."
}

# targetIndex now refers to the paragraph containing "{% endif %}" that
# closes the Skills section. Insert the 4 new paragraphs right after it.

$r = $d.Paragraphs.Item($targetIndex).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Item($targetIndex + 1).Range.InsertAfter("{% if languages %}")

$r = $d.Paragraphs.Item($targetIndex + 1).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$langHeadingPara = $d.Paragraphs.Item($targetIndex + 2)
$langHeadingPara.Style = $d.Styles.Item("Heading 1")
$langHeadingPara.Range.InsertAfter("Languages")

$r = $d.Paragraphs.Item($targetIndex + 2).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$langForPara = $d.Paragraphs.Item($targetIndex + 3)
$langForPara.Style = $d.Styles.Item("Normal")
$langForPara.Range.InsertAfter("{% for lang in languages %}{{ lang }}{% if not loop.last %}, {% endif %}{% endfor %}")

$r = $d.Paragraphs.Item($targetIndex + 3).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$closingEndifPara = $d.Paragraphs.Item($targetIndex + 4)
$closingEndifPara.Style = $d.Styles.Item("Normal")
$closingEndifPara.Range.InsertAfter("{% endif %}")

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
